$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 105 (shifts old 105-107 down to 107-109)
$ws.Rows.Item(105).Insert()
$ws.Rows.Item(105).Insert()

# New row 105: weekly update for Region Metropolitana, $/caja 18 unidades
$ws.Cells.Item(105, 1).Value = 10
$ws.Cells.Item(105, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(105, 3).Value = "La Araucanía"
$ws.Cells.Item(105, 4).Value = 45075
$ws.Cells.Item(105, 5).Value = 9
$ws.Cells.Item(105, 6).Value = 100112010
$ws.Cells.Item(105, 7).Value = "Achicoria"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 65
$ws.Cells.Item(105, 11).Value = 10000
$ws.Cells.Item(105, 12).Value = 10000
$ws.Cells.Item(105, 13).Value = 10000
$ws.Cells.Item(105, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(105, 15).Value = "Región Metropolitana"
$ws.Cells.Item(105, 16).Value = 556
$ws.Cells.Item(105, 17).Value = 18
$ws.Cells.Item(105, 18).Value = "Hortaliza"

# New row 106: weekly update for Region del Maule, $/caja 18 unidades
$ws.Cells.Item(106, 1).Value = 10
$ws.Cells.Item(106, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(106, 3).Value = "La Araucanía"
$ws.Cells.Item(106, 4).Value = 45075
$ws.Cells.Item(106, 5).Value = 9
$ws.Cells.Item(106, 6).Value = 100112010
$ws.Cells.Item(106, 7).Value = "Achicoria"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 65
$ws.Cells.Item(106, 11).Value = 7500
$ws.Cells.Item(106, 12).Value = 7500
$ws.Cells.Item(106, 13).Value = 7500
$ws.Cells.Item(106, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(106, 15).Value = "Región del Maule"
$ws.Cells.Item(106, 16).Value = 417
$ws.Cells.Item(106, 17).Value = 18
$ws.Cells.Item(106, 18).Value = "Hortaliza"
